{"js": "// Recolor specific table-cell shading from dark gray (#A9A9A9) to white (#FFFFFF).\n// The document contains two tables (each 13 rows x 4 cols). A fixed subset of the\n// dark-gray shaded data cells is updated to white; the remaining dark-gray cells\n// (which encode separate, unrelated data) are left untouched.\n\n// (tableIndex, rowIndex, columnIndex) \u2014 all 0-based, matching document order.\nconst targets = [\n  [0, 1, 3], [0, 3, 3], [0, 4, 2], [0, 5, 3],\n  [0, 8, 1], [0, 8, 3], [0, 9, 2], [0, 10, 3],\n  [0, 11, 2], [0, 12, 1], [0, 12, 3],\n  [1, 1, 1], [1, 1, 3], [1, 3, 1], [1, 3, 3],\n  [1, 4, 1], [1, 4, 2], [1, 4, 3], [1, 5, 3],\n  [1, 6, 3], [1, 7, 2], [1, 8, 3], [1, 9, 1],\n  [1, 9, 2], [1, 9, 3], [1, 10, 3], [1, 11, 2],\n  [1, 11, 3], [1, 12, 3]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const [tableIndex, rowIndex, columnIndex] of targets) {\n  const table = tables.items[tableIndex];\n  const cell = table.getCell(rowIndex, columnIndex);\n  cell.shadingColor = \"#FFFFFF\";\n}\n\nawait context.sync();\n", "ps1": "# Recolor specific table-cell shading from dark gray (#A9A9A9) to white (#FFFFFF).\n# The document contains two tables (each 13 rows x 4 cols). A fixed subset of the\n# dark-gray shaded data cells is updated to white; the remaining dark-gray cells\n# (which encode separate, unrelated data) are left untouched.\n\n$d = $word.ActiveDocument\n\n# (tableIndex, rowIndex, columnIndex) using 1-based COM indices.\n$targets = @(\n    @(1, 2, 4), @(1, 4, 4), @(1, 5, 3), @(1, 6, 4),\n    @(1, 9, 2), @(1, 9, 4), @(1, 10, 3), @(1, 11, 4),\n    @(1, 12, 3), @(1, 13, 2), @(1, 13, 4),\n    @(2, 2, 2), @(2, 2, 4), @(2, 4, 2), @(2, 4, 4),\n    @(2, 5, 2), @(2, 5, 3), @(2, 5, 4), @(2, 6, 4),\n    @(2, 7, 4), @(2, 8, 3), @(2, 9, 4), @(2, 10, 2),\n    @(2, 10, 3), @(2, 10, 4), @(2, 11, 4), @(2, 12, 3),\n    @(2, 12, 4), @(2, 13, 4)\n)\n\nforeach ($target in $targets) {\n    $tableIndex = $target[0]\n    $rowIndex = $target[1]\n    $columnIndex = $target[2]\n    $tbl = $d.Tables.Item($tableIndex)\n    $cell = $tbl.Cell($rowIndex, $columnIndex)\n    $cell.Shading.BackgroundPatternColor = 16777215\n}\n"}
